$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.432.43'
$ws.Range('E2').Value = '  +0.96%  '

$ws.Range('D3').Value = '2.614.37'
$ws.Range('E3').Value = '  +9.74%  '

$ws.Range('E4').Value = '  -0.18%  '

$ws.Range('D5').Value = '''307.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.02%  '

$ws.Range('D6').Value = '''101.15'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.32%  '

$ws.Range('D7').Value = '''0.605'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +6.16%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('E9').Value = '  +12.60%  '

$ws.Range('D10').Value = '''39.10'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +13.63%  '

$ws.Range('E11').Value = '  +6.13%  '

$ws.Range('D12').Value = '''8.20'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +14.52%  '

$ws.Range('D13').Value = '3.012.96'
$ws.Range('E13').Value = '  +9.56%  '

$ws.Range('E14').Value = '  +1.96%  '

$ws.Range('D15').Value = '2.622.15'
$ws.Range('E15').Value = '  +10.03%  '

$ws.Range('E16').Value = '  +10.92%  '

$ws.Range('D17').Value = '''14.92'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +9.01%  '

$ws.Range('D18').Value = '46.609.20'
$ws.Range('E18').Value = '  +1.42%  '

$ws.Range('D19').Value = '''13.35'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.43%  '

$ws.Range('E20').Value = '  +5.61%  '

$ws.Range('D21').Value = '''6.71'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +11.15%  '

$ws.Range('D22').Value = '''71.16'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.33%  '

$ws.Range('D23').Value = '''257.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.11%  '

$ws.Range('E24').Value = '  +7.26%  '

$ws.Range('E25').Value = '  +16.11%  '

$ws.Range('D26').Value = '''28.43'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +34.54%  '

$ws.Range('D27').Value = '''1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.01%  '

$ws.Range('E28').Value = '  +7.42%  '

$ws.Range('D29').Value = '''40.44'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.18%  '

$ws.Range('E30').Value = '  +3.44%  '

$ws.Range('E31').Value = '  +12.12%  '

$ws.Range('D32').Value = '''3.75'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.82%  '

$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').Value = '''2.34'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +20.67%  '

$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = '''2.96'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.22%  '

$ws.Range('D35').Value = '''0.0836'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.60%  '

$ws.Range('D36').Value = '''151.01'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.78%  '

$ws.Range('E37').Value = '  +4.55%  '

$ws.Range('E38').Value = '  +5.36%  '

$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '''4.19'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.66%  '

$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = '''15.82'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.47%  '

$ws.Range('E41').Value = '  +13.20%  '

$ws.Range('E42').Value = '  +7.95%  '

$ws.Range('D43').Value = '2.047.68'
$ws.Range('E43').Value = '  +6.41%  '

$ws.Range('D44').Value = '''19.54'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +38.14%  '

$ws.Range('E45').Value = '  -0.10%  '

$ws.Range('D46').Value = '''91.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.31%  '

$ws.Range('D47').Value = '''9.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +10.19%  '

$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '''110.40'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +12.54%  '

$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '''1.79'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.12%  '

$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.870.97'
$ws.Range('E50').Value = '  +9.41%  '

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '''0.201'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.75%  '
